$p = $ppt.ActivePresentation

# --- Chart 1 (Server Side Frameworks Popularity) on slide 3 ---
$slide3 = $p.Slides.Item(3)
$chart1 = $slide3.Shapes.Item(2).Chart
$legend1 = $chart1.Legend

# Enlarge the legend font and give it more room (matches the wider/taller
# manual layout box the legend needs once the font is bigger).
$legend1.Left   = 177.7487401574803
$legend1.Top    = 8.97503937007874
$legend1.Width  = 543.5024409448819
$legend1.Height = 48.08291338582678
$legend1.Font.Size = 20

# --- Chart 2 (Client Side Frameworks Popularity) on slide 5 ---
$slide5 = $p.Slides.Item(5)
$chart2 = $slide5.Shapes.Item(2).Chart
$legend2 = $chart2.Legend

$legend2.Left   = 216.0
$legend2.Top    = 8.97503937007874
$legend2.Width  = 320.5890551181102
$legend2.Height = 88.08291338582677
$legend2.Font.Size = 18
